# Quina_edt.xlsx - append the newest 12 draws (concursos 6940-6951) to the
# bottom of the QUINA results table on sheet "QUINA".
#
# Before: data runs from row 2 through row 448 (dimension A1:F448).
# After:  12 new rows (449-460) are appended, each with the contest number
#         in column A and the five drawn numbers in columns B-F.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Concurso, Bola1, Bola2, Bola3, Bola4, Bola5
$novosResultados = @(
    @(6940, 24, 53, 66, 73, 77),
    @(6941, 12, 32, 34, 57, 64),
    @(6942, 16, 33, 34, 50, 71),
    @(6943, 22, 23, 35, 40, 44),
    @(6944,  2,  8, 30, 56, 61),
    @(6945, 33, 61, 66, 68, 70),
    @(6946,  1, 48, 53, 75, 80),
    @(6947,  6, 30, 52, 60, 79),
    @(6948,  3, 21, 32, 46, 57),
    @(6949, 21, 51, 60, 67, 73),
    @(6950,  1,  6, 24, 47, 60),
    @(6951,  1, 10, 20, 44, 66)
)

$primeiraLinha = 449
$linha = $primeiraLinha
foreach ($sorteio in $novosResultados) {
    for ($col = 1; $col -le 6; $col++) {
        $ws.Cells.Item($linha, $col).Value = $sorteio[$col - 1]
    }
    $linha++
}
$ultimaLinha = $linha - 1

# Reflect where the user ended up working: anchored on the first cell of the
# freshly typed last row, with the rest of the newly entered block selected.
$ws.Range("C" + $primeiraLinha + ":F" + $ultimaLinha).Select() | Out-Null
